# LoginTestData.xlsx - "Configuration issue resolved and all the tests are
# working as expected"
#
# The only real content change is on the "AddAlbum" sheet: the test data
# value in A2 was renamed from "Fest2" to "FestT" (fixing a typo that was
# causing the configuration/test issue). Saving that change through Excel
# also re-numbers the shared-strings table and records the new selection
# (A2) on that sheet, while leaving every other sheet's active selection
# and the workbook's active tab (ViewMonthlyAttendance) untouched.

$wb = $excel.ActiveWorkbook

$wsAddAlbum = $wb.Worksheets.Item("AddAlbum")
$wsViewMonthlyAttendance = $wb.Worksheets.Item("ViewMonthlyAttendance")

# Fix the test data value.
$wsAddAlbum.Range("A2").Value = "FestT"

# Record the new selection on the AddAlbum sheet (A2) without leaving it as
# the active sheet.
$wsAddAlbum.Activate() | Out-Null
$wsAddAlbum.Range("A2").Select() | Out-Null

# Restore the originally active sheet/tab.
$wsViewMonthlyAttendance.Activate() | Out-Null
